# Updated symbol list on Wed Dec 14 12:30:52 UTC 2022 with GitHub Actions
#
# All data cells on Sheet1 are stored as text (the "Price" column D and the
# "Hora" column G look numeric but must stay text), so every numeric-looking
# replacement is entered with a leading apostrophe to force Excel to treat it
# as a literal string instead of auto-converting it to a number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - BNB
$ws.Range("D2").Value = "'273.18"
$ws.Range("G2").Value = "'12"

# Row 3 - OKB
$ws.Range("D3").Value = "'22.92"
$ws.Range("G3").Value = "'12"

# Row 4 - HuobiToken
$ws.Range("D4").Value = "'6.368"
$ws.Range("G4").Value = "'12"

# Row 5 - Cronos
$ws.Range("D5").Value = "'0.06217"
$ws.Range("G5").Value = "'12"

# Row 6 - GateToken
$ws.Range("D6").Value = "'3.653"
$ws.Range("G6").Value = "'12"

# Row 7 - KuCoinToken
$ws.Range("D7").Value = "'6.709"
$ws.Range("G7").Value = "'12"

# Row 8 - FTXToken
$ws.Range("D8").Value = "'1.381"
$ws.Range("G8").Value = "'12"

# Row 9 - MXToken (only Hora changes)
$ws.Range("G9").Value = "'12"

# Row 10 - One
$ws.Range("D10").Value = "'0.01376"
$ws.Range("G10").Value = "'12"

# Row 11 - WazirX
$ws.Range("D11").Value = "'0.1632"
$ws.Range("G11").Value = "'12"

# Row 12 - MandalaExchangeToken
$ws.Range("D12").Value = "'0.08331"
$ws.Range("G12").Value = "'12"

# Row 13 - LiechtensteinCryptoassetsExchange
$ws.Range("D13").Value = "'0.03394"
$ws.Range("G13").Value = "'12"

# Row 14 - BitrueCoin
$ws.Range("D14").Value = "'0.03104"
$ws.Range("G14").Value = "'12"

# Row 15
$ws.Range("D15").Value = "'0.09332"
$ws.Range("G15").Value = "'12"

# Row 16
$ws.Range("D16").Value = "'3.871"
$ws.Range("G16").Value = "'12"

# Row 17
$ws.Range("D17").Value = "'0.001639"
$ws.Range("G17").Value = "'12"

# Row 18 - CoinExToken
$ws.Range("D18").Value = "'0.04811"
$ws.Range("G18").Value = "'12"

# Row 19 - TigerCash
$ws.Range("D19").Value = "'0.006221"
$ws.Range("G19").Value = "'12"

# Row 20 and 21 swap places: HotbitToken <-> BitKan
$ws.Range("B20").Value = "BitKan"
$ws.Range("C20").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("D20").Value = "'0.001087"
$ws.Range("E20").Value = "19BitKanKAN"
$ws.Range("G20").Value = "'12"

$ws.Range("B21").Value = "HotbitToken"
$ws.Range("C21").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("D21").Value = "'0.003191"
$ws.Range("E21").Value = "20HotbitTokenHTBWorstin24h"
$ws.Range("G21").Value = "'12"

# Row 22 - NitroEx
$ws.Range("D22").Value = "'0.0001499"
$ws.Range("G22").Value = "'12"

# Row 23
$ws.Range("D23").Value = "'3.732"
$ws.Range("G23").Value = "'12"

# Row 24
$ws.Range("D24").Value = "'2.379"
$ws.Range("G24").Value = "'12"

# Row 25
$ws.Range("D25").Value = "'0.3402"
$ws.Range("G25").Value = "'12"

# Rows 26-39 - only the Hora column changes
$ws.Range("G26").Value = "'12"
$ws.Range("G27").Value = "'12"
$ws.Range("G28").Value = "'12"
$ws.Range("G29").Value = "'12"
$ws.Range("G30").Value = "'12"
$ws.Range("G31").Value = "'12"
$ws.Range("G32").Value = "'12"
$ws.Range("G33").Value = "'12"
$ws.Range("G34").Value = "'12"
$ws.Range("G35").Value = "'12"
$ws.Range("G36").Value = "'12"
$ws.Range("G37").Value = "'12"
$ws.Range("G38").Value = "'12"
$ws.Range("G39").Value = "'12"

# Row 40 - IDEX
$ws.Range("D40").Value = "'0.04685"
$ws.Range("G40").Value = "'12"

# Row 41 - KickToken
$ws.Range("D41").Value = "'0.006931"
$ws.Range("G41").Value = "'12"

# Row 42 - BKEXToken
$ws.Range("D42").Value = "'0.1170"
$ws.Range("G42").Value = "'12"

# Row 43 - CEJI
$ws.Range("D43").Value = "'0.003448"
$ws.Range("E43").Value = "42CEJICEJI"
$ws.Range("G43").Value = "'12"

# Row 44 - LocalTraders
$ws.Range("D44").Value = "'0.01159"
$ws.Range("G44").Value = "'12"

# Row 45 - CoinLion
$ws.Range("D45").Value = "'0.00006249"
$ws.Range("G45").Value = "'12"

# Row 46 - Kangarootoken (only Hora changes)
$ws.Range("G46").Value = "'12"

# Row 47 - CoinbaseStockToken
$ws.Range("D47").Value = "'0.8994"
$ws.Range("G47").Value = "'12"

# Row 48 - BOLO
$ws.Range("D48").Value = "'0.05173"
$ws.Range("G48").Value = "'12"

# Row 49 - CryptobidCoin
$ws.Range("D49").Value = "'0.00002299"
$ws.Range("G49").Value = "'12"

# Row 50 - SpecialPowerGold
$ws.Range("D50").Value = "'0.01239"
$ws.Range("G50").Value = "'12"

# Row 51 - DigiFinexToken (only Hora changes)
$ws.Range("G51").Value = "'12"
